$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting (style) from row 14 to the two new rows 15 and 16
$ws.Range("A14:J14").Copy() | Out-Null
$ws.Range("A15:J15").PasteSpecial(-4122) | Out-Null
$ws.Range("A16:J16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Populate data rows 2-16 with the refreshed convergence results
# Row 2
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = "Poisson"
$ws.Range("C2").Value = "FE"
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = "Regular_RightTriangles"
$ws.Range("F2").Value = "Dirichlet"
$ws.Range("G2").Value = 2.0039
$ws.Range("H2").Value = "Triangles"
$ws.Range("I2").Value = "Green"
$ws.Range("J2").Value = 103.069

# Row 3
$ws.Range("A3").Value = 10
$ws.Range("B3").Value = "Poisson"
$ws.Range("C3").Value = "FE"
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = "Regular_RightTriangles"
$ws.Range("F3").Value = "Neumann"
$ws.Range("G3").Value = 0.9103
$ws.Range("H3").Value = "Triangles"
$ws.Range("I3").Value = "Orange (ILU)"
$ws.Range("J3").Value = 26.158

# Row 4
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Poisson"
$ws.Range("C4").Value = "FE"
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = "Unstructured_triangles"
$ws.Range("F4").Value = "Dirichlet"
$ws.Range("G4").Value = 2.0156
$ws.Range("H4").Value = "Triangles"
$ws.Range("I4").Value = "Green"
$ws.Range("J4").Value = 7.609

# Row 5
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Poisson"
$ws.Range("C5").Value = "FE"
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = "Unstructured_triangles"
$ws.Range("F5").Value = "Neumann"
$ws.Range("G5").Value = 0.8202
$ws.Range("H5").Value = "Squares"
$ws.Range("I5").Value = "Red"
$ws.Range("J5").Value = 3.102

# Row 6
$ws.Range("A6").Value = 7
$ws.Range("B6").Value = "Poisson"
$ws.Range("C6").Value = "FE"
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = "Regular_Tetrahedra"
$ws.Range("F6").Value = "Dirichlet"
$ws.Range("G6").Value = 1.3403
$ws.Range("H6").Value = "Tetrahedron"
$ws.Range("I6").Value = "Green"
$ws.Range("J6").Value = 208.494

# Row 7
$ws.Range("A7").Value = 8
$ws.Range("B7").Value = "Poisson"
$ws.Range("C7").Value = "FE"
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = "Unstructured_Tetrahedra"
$ws.Range("F7").Value = "Dirichlet"
$ws.Range("G7").Value = 0.6691
$ws.Range("H7").Value = "Tetrahedron"
$ws.Range("I7").Value = "Green"
$ws.Range("J7").Value = 12.001

# Row 8
$ws.Range("A8").Value = 1
$ws.Range("B8").Value = "Poisson"
$ws.Range("C8").Value = "FV"
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = "RegularSquares"
$ws.Range("F8").Value = "Dirichlet"
$ws.Range("G8").Value = 2.0039
$ws.Range("H8").Value = "Squares"
$ws.Range("I8").Value = "Green"
$ws.Range("J8").Value = 9.776999999999999

# Row 9
$ws.Range("A9").Value = 13
$ws.Range("B9").Value = "Poisson"
$ws.Range("C9").Value = "FV"
$ws.Range("D9").Value = 2
$ws.Range("E9").Value = "RegularSquares"
$ws.Range("F9").Value = "Neumann"
$ws.Range("G9").Value = 2.0039
$ws.Range("H9").Value = "Squares"
$ws.Range("I9").Value = "Green"
$ws.Range("J9").Value = 9.827

# Row 10
$ws.Range("A10").Value = 5
$ws.Range("B10").Value = "Poisson"
$ws.Range("C10").Value = "FV"
$ws.Range("D10").Value = 2
$ws.Range("E10").Value = "Regular_RightTriangles"
$ws.Range("F10").Value = "Dirichlet"
$ws.Range("G10").Value = 0.0212
$ws.Range("H10").Value = "Triangles"
$ws.Range("I10").Value = "Green"
$ws.Range("J10").Value = 15.552

# Row 11
$ws.Range("A11").Value = 0
$ws.Range("B11").Value = "Poisson"
$ws.Range("C11").Value = "FV"
$ws.Range("D11").Value = 2
$ws.Range("E11").Value = "Regular_RightTriangles"
$ws.Range("F11").Value = "Neumann"
$ws.Range("G11").Value = -0.0056
$ws.Range("H11").Value = "Triangles"
$ws.Range("I11").Value = "Orange `n (suspicious order 0 convergence)"
$ws.Range("J11").Value = 16.075

# Row 12
$ws.Range("A12").Value = 6
$ws.Range("B12").Value = "Poisson"
$ws.Range("C12").Value = "FV"
$ws.Range("D12").Value = 2
$ws.Range("E12").Value = "Structured_triangles"
$ws.Range("F12").Value = "Dirichlet"
$ws.Range("G12").Value = 0.8952
$ws.Range("H12").Value = "Triangles"
$ws.Range("I12").Value = "Orange `n (BC don't fit the domain)"
$ws.Range("J12").Value = 4.925

# Row 13
$ws.Range("A13").Value = 9
$ws.Range("B13").Value = "Poisson"
$ws.Range("C13").Value = "FV"
$ws.Range("D13").Value = 2
$ws.Range("E13").Value = "Unstructured_triangles"
$ws.Range("F13").Value = "Dirichlet"
$ws.Range("G13").Value = 0.6138
$ws.Range("H13").Value = "Triangles"
$ws.Range("I13").Value = "Green"
$ws.Range("J13").Value = 2.235

# Row 14
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "Poisson"
$ws.Range("C14").Value = "FV"
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = "Regular_Cubes"
$ws.Range("F14").Value = "Dirichlet"
$ws.Range("G14").Value = 1.3403
$ws.Range("H14").Value = "Cubes"
$ws.Range("I14").Value = "Green"
$ws.Range("J14").Value = 5.891

# Row 15
$ws.Range("A15").Value = 11
$ws.Range("B15").Value = "Poisson"
$ws.Range("C15").Value = "FV"
$ws.Range("D15").Value = 3
$ws.Range("E15").Value = "Regular_Tetrahedra"
$ws.Range("F15").Value = "Dirichlet"
$ws.Range("G15").Value = 0.0065
$ws.Range("H15").Value = "Tetrahedron"
$ws.Range("I15").Value = "Green"
$ws.Range("J15").Value = 62.301

# Row 16
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Poisson"
$ws.Range("C16").Value = "FV"
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = "Unstructured_Tetrahedra"
$ws.Range("F16").Value = "Dirichlet"
$ws.Range("G16").Value = 0.5359
$ws.Range("H16").Value = "Tetrahedron"
$ws.Range("I16").Value = "Green"
$ws.Range("J16").Value = 3.712

# Re-fit rows containing multi-line text so no explicit custom row height is stored
$ws.Rows.Item(11).AutoFit() | Out-Null
$ws.Rows.Item(12).AutoFit() | Out-Null
